# Update "paises.xlsx" (sheet "Pais") with the latest COVID-19 country
# statistics snapshot and refresh the "last updated" timestamp.
#
# Columns: A=Pais  B=Casos totales  C=Nuevos casos  D=Casos activos
#          E=Recuperados  F=Casos criticos  G=Muertes hoy  H=Muertes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (A1) ------------------------------------------------
$ws.Range("A1").Value2 = "Datos actualizados a 19 de Octubre de 2020 a las 20:49"

# --- Row 4: Estados Unidos -------------------------------------------------
$ws.Range("B4").Value2 = 8408992
$ws.Range("C4").Value2 = 21193
$ws.Range("D4").Value2 = 5472482
$ws.Range("E4").Value2 = 2711562
$ws.Range("G4").Value2 = 218
$ws.Range("H4").Value2 = 224948

# --- Row 5: India ------------------------------------------------------
$ws.Range("B5").Value2 = 7593994
$ws.Range("C5").Value2 = 45756
$ws.Range("D5").Value2 = 6730054
$ws.Range("E5").Value2 = 748707
$ws.Range("G5").Value2 = 591
$ws.Range("H5").Value2 = 115233

# --- Row 11: Francia ----------------------------------------------------
$ws.Range("B11").Value2 = 910277
$ws.Range("C11").Value2 = 13243
$ws.Range("D11").Value2 = 105935
$ws.Range("E11").Value2 = 770719
$ws.Range("G11").Value2 = 146
$ws.Range("H11").Value2 = 33623

# --- Row 21: Alemania -----------------------------------------------------
$ws.Range("B21").Value2 = 372255
$ws.Range("C21").Value2 = 5274
$ws.Range("E21").Value2 = 70463
$ws.Range("G21").Value2 = 26
$ws.Range("H21").Value2 = 9892

# --- Row 24: Arabia Saudita --------------------------------------------------
$ws.Range("B24").Value2 = 349451
$ws.Range("C24").Value2 = 1958
$ws.Range("E24").Value2 = 34653

# --- Row 31: Canada --------------------------------------------------
$ws.Range("B31").Value2 = 199970
$ws.Range("C31").Value2 = 1822
$ws.Range("D31").Value2 = 168699
$ws.Range("E31").Value2 = 21499
$ws.Range("G31").Value2 = 12
$ws.Range("H31").Value2 = 9772

# --- Row 36: Marruecos -------------------------------------------------
$ws.Range("B36").Value2 = 153423
$ws.Range("C36").Value2 = 134
$ws.Range("E36").Value2 = 6841
$ws.Range("G36").Value2 = 8
$ws.Range("H36").Value2 = 12395

# --- Row 52: Etiopia -------------------------------------------------
$ws.Range("B52").Value2 = 89860
$ws.Range("C52").Value2 = 723
$ws.Range("D52").Value2 = 43149
$ws.Range("E52").Value2 = 45346
$ws.Range("G52").Value2 = 13
$ws.Range("H52").Value2 = 1365

# --- Row 62: Uzbekistan -------------------------------------------------
$ws.Range("B62").Value2 = 63523
$ws.Range("C62").Value2 = 399
$ws.Range("D62").Value2 = 60604
$ws.Range("E62").Value2 = 2388
$ws.Range("G62").Value2 = 6
$ws.Range("H62").Value2 = 531

# --- Row 67: Argelia ----------------------------------------------------
$ws.Range("B67").Value2 = 54616
$ws.Range("C67").Value2 = 214
$ws.Range("D67").Value2 = 38215
$ws.Range("E67").Value2 = 14536
$ws.Range("G67").Value2 = 9
$ws.Range("H67").Value2 = 1865

# --- Row 93: Costa de Marfil ------------------------------------------------
$ws.Range("B93").Value2 = 20324
$ws.Range("C93").Value2 = 1
$ws.Range("D93").Value2 = 20029
$ws.Range("E93").Value2 = 174

# --- Row 106: Maldivas (unchanged position, refreshed stats) --------------
$ws.Range("B106").Value2 = 11232
$ws.Range("C106").Value2 = 22
$ws.Range("D106").Value2 = 10201
$ws.Range("E106").Value2 = 994

# --- Row 107: now Mozambique (moved up, overtaking the next two rows) -----
$ws.Range("A107").Value2 = "Mozambique"
$ws.Range("B107").Value2 = 11080
$ws.Range("C107").Value2 = 214
$ws.Range("D107").Value2 = 8836
$ws.Range("E107").Value2 = 2169
$ws.Range("G107").Value2 = 0
$ws.Range("H107").Value2 = 75

# --- Row 108: now Consejo Danes para los Refugiados (pushed down one) -----
$ws.Range("A108").Value2 = "Consejo Danes para los Refugiados"
$ws.Range("B108").Value2 = 11052
$ws.Range("C108").Value2 = 46
$ws.Range("D108").Value2 = 10357
$ws.Range("E108").Value2 = 392
$ws.Range("G108").Value2 = 1
$ws.Range("H108").Value2 = 303

# --- Row 109: now Luxemburgo (pushed down one) -----------------------------
$ws.Range("A109").Value2 = "Luxemburgo"
$ws.Range("B109").Value2 = 11010
$ws.Range("C109").Value2 = 122
$ws.Range("D109").Value2 = 8471
$ws.Range("E109").Value2 = 2404
$ws.Range("G109").Value2 = 2
$ws.Range("H109").Value2 = 135

# Row 110 (Uganda) is unchanged.

# --- Row 116: Zimbabue -------------------------------------------------------
$ws.Range("B116").Value2 = 8159
$ws.Range("C116").Value2 = 12
$ws.Range("D116").Value2 = 7683
$ws.Range("E116").Value2 = 244
$ws.Range("G116").Value2 = 1
$ws.Range("H116").Value2 = 232

# --- Row 124: Suazilandia ------------------------------------------------
$ws.Range("B124").Value2 = 5788
$ws.Range("C124").Value2 = 8
$ws.Range("D124").Value2 = 5427
$ws.Range("E124").Value2 = 245

# --- Row 127: Republica de Yibuti ------------------------------------------------
$ws.Range("B127").Value2 = 5469
$ws.Range("C127").Value2 = 10
$ws.Range("D127").Value2 = 5379
$ws.Range("E127").Value2 = 29

# --- Row 154: Republica de Chipre ------------------------------------------------
$ws.Range("B154").Value2 = 2687
$ws.Range("C154").Value2 = 43
$ws.Range("E154").Value2 = 1218

# --- Row 165: Republica del Chad ------------------------------------------------
$ws.Range("B165").Value2 = 1390
$ws.Range("C165").Value2 = 11
$ws.Range("D165").Value2 = 1194
$ws.Range("E165").Value2 = 103
